# Apply the setAlignmentNotExisting stack-trace text update:
#  - MParagraphImpl hash code changed (@7da34b26 -> @51b77cdf)
#  - various stack frame line numbers / frame counts shifted
#    (library version bump referenced in the commit message)
$d = $word.ActiveDocument
$wdReplaceAll = 2
$wdFindContinue = 1

$old0 = "MParagraphImpl@7da34b26"
$new0 = "MParagraphImpl@51b77cdf"
$r0 = $d.Content.Find.Execute($old0, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new0, $wdReplaceAll)
Write-Output "op0: $r0"

$old1 = "`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)"
$new1 = "`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:192)"
$r1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new1, $wdReplaceAll)
Write-Output "op1: $r1"

$old2 = "`tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:109)"
$new2 = "`tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)"
$r2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new2, $wdReplaceAll)
Write-Output "op2: $r2"

$old3 = "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:591)"
$new3 = "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:604)"
$r3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new3, $wdReplaceAll)
Write-Output "op3: $r3"

$old4 = "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1331)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1556)"
$new4 = "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1675)"
$r4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new4, $wdReplaceAll)
Write-Output "op4: $r4"

$old5 = "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1331)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:301)"
$new5 = "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:314)"
$r5 = $d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new5, $wdReplaceAll)
Write-Output "op5: $r5"

$old6 = "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1331)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:286)"
$new6 = "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:299)"
$r6 = $d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new6, $wdReplaceAll)
Write-Output "op6: $r6"

$old7 = "`tat sun.reflect.GeneratedMethodAccessor4.invoke(Unknown Source)"
$new7 = "`tat sun.reflect.GeneratedMethodAccessor6.invoke(Unknown Source)"
$r7 = $d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new7, $wdReplaceAll)
Write-Output "op7: $r7"

$old8 = "`tat org.junit.runners.model.FrameworkMethod`$1.runReflectiveCall(FrameworkMethod.java:59)"
$new8 = "`tat org.junit.runners.model.FrameworkMethod`$1.runReflectiveCall(FrameworkMethod.java:50)"
$r8 = $d.Content.Find.Execute($old8, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new8, $wdReplaceAll)
Write-Output "op8: $r8"

$old9 = "`tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:56)"
$new9 = "`tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)"
$r9 = $d.Content.Find.Execute($old9, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new9, $wdReplaceAll)
Write-Output "op9: $r9"

$old10 = "`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.BlockJUnit4ClassRunner`$1.evaluate(BlockJUnit4ClassRunner.java:100)`n`tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:366)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:103)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:63)`n`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)"
$new10 = "`tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)`n`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)"
$r10 = $d.Content.Find.Execute($old10, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new10, $wdReplaceAll)
Write-Output "op10: $r10"

$old11 = "`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)"
$new11 = "`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)"
$r11 = $d.Content.Find.Execute($old11, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new11, $wdReplaceAll)
Write-Output "op11: $r11"

$old12 = "`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)"
$new12 = "`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)"
$r12 = $d.Content.Find.Execute($old12, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new12, $wdReplaceAll)
Write-Output "op12: $r12"

$old13 = "`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)"
$new13 = "`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)"
$r13 = $d.Content.Find.Execute($old13, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new13, $wdReplaceAll)
Write-Output "op13: $r13"
